$wb = $excel.ActiveWorkbook

# --- Update Metadata sheet timestamp (A2) ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 09:59 AM"

# --- Update Stock List sheet: new top row inserted, existing rows shifted down ---
$ws = $wb.Worksheets.Item("Stock List")

$data = @(
    @("CAPTRU-RE1", 5.67, -11.9565, 0),
    @("NIFTYCASE", 10.19, -0.5854, 0),
    @("MOMENTUM30", 31.54, -0.6614, 0),
    @("CANHLIFE", 118.46, 0.6286, 11253.7),
    @("FLEXIADD", 10.64, -1.0233, 0),
    @("MOENERGY", 36.3, -0.6568000000000001, 0),
    @("MONIFTY100", 26.49, 0.3409, 0),
    @("RUBICON", 652.65, -0.1453, 10752.4289),
    @("CRAMC", 317.2, 2.3226, 6325.5208),
    @("LGEINDIA", 1633.4, -0.946, 110870.6825),
    @("TATACAP", 329.3, 0.1521, 139783.5374),
    @("ELIQUID", 1004.85, 0.0408, 0),
    @("WEWORK", 632.15, -2.4008, 8472.2803),
    @("GROWWRLTY", 10.8, -0.4608, 0),
    @("ADVANCE", 130.05, -5.2666, 836.0358),
    @("OMFREIGHT", 88.90000000000001, -0.5926, 299.3747),
    @("GLOTTIS", 72.73999999999999, -0.8587, 672.1394),
    @("FABTECH", 237.72, 0.4734, 1056.6843),
    @("PACEDIGITK", 218.85, 0.1327, 4723.9063),
    @("JAINREC", 377.25, 1.2208, 13018.3623),
    @("EPACKPEB", 301.45, 1.979, 3028.1254),
    @("BMWVENTLTD", 69.25, 0, 600.5014),
    @("STYL", 372.4, -0.8388, 6025.649),
    @("JARO", 621.5, -1.4821, 1377.0134),
    @("SOLARWORLD", 309.1, -0.6269, 2679.0517),
    @("ARSSBL", 537.3, 4.7266, 3370.2277),
    @("GANESHCP", 274.4, -2.7984, 1108.9312),
    @("ATLANTAELE", 1003.05, -1.7436, 7713.116),
    @("GKENERGY", 213.85, -0.7933, 4337.2472),
    @("SAATVIKGL", 528.2, -1.3079, 6713.6863),
    @("IVALUE", 281.45, -0.3364, 1506.8799),
    @("VMSTMT", 70.03, -0.9056, 347.5674),
    @("EUROPRATIK", 321.75, 0.8147, 3288.285),
    @("SHRINGARMS", 229.31, -1.2616, 2211.284),
    @("DEVX", 44.53, -0.3803, 401.605),
    @("URBANCO", 148.9, -2.0459, 21380.5798),
    @("SML100CASE", 10.36, -0.7663, 0),
    @("AONEGOLD", 11.28, -0.2653, 0),
    @("ELM250", 16.72, 0.1797, 0),
    @("AMANTA", 122.52, 1.407, 475.7372),
    @("CPEDU", 315.9, 1.8539, 574.7148999999999),
    @("AHCL", 139.27, 3.1706, 740.2409),
    @("STLNETWORK", 26.59, -0.412, 1297.3822),
    @("VIKRAN", 98.05, -1.783, 2528.8166),
    @("MANUFGBEES", 151.77, -1.011, 0),
    @("MEIL", 461.15, -0.7319, 1274.1632),
    @("GROWWNXT50", 70.29000000000001, -0.4109, 0),
    @("SHREEJISPG", 270.05, -0.7899, 4399.6074),
    @("GEMAROMA", 219.52, -0.876, 1146.7097),
    @("PATELRMART", 219.31, -1.0646, 732.5069999999999),
    @("VIKRAMSOLR", 322, -1.5892, 11647.2884),
    @("LTGILTCASE", 29.67, 0.2365, 0),
    @("REGAAL", 89.13, -0.8675, 915.5742),
    @("BLUESTONE", 711.95, 0.1266, 10773.2539),
    @("MOSILVER", 145.9, -1.5054, 0),
    @("ALLTIME", 308.75, 2.66, 2022.5526),
    @("JSWCEMENT", 134.98, -0.4793, 18402.6999),
    @("SBILIQETF", 1012.94, 0.0296, 0),
    @("HILINFRA", 77.23, -0.3998, 0),
    @("GROWWPOWER", 10.28, -0.9634, 0),
    @("LOTUSDEV", 177.82, 0.3669, 8690.485000000001),
    @("MBEL", 450.2, -0.7714, 2572.8126),
    @("LAXMIINDIA", 145.62, -1.1942, 761.1248000000001),
    @("CPPLUS", 1322.1, -0.264, 15497.9053),
    @("SHANTIGOLD", 241.57, -1.6409, 1741.6231),
    @("MOGOLD", 119.65, -0.5403, 0),
    @("BRIGHOTEL", 82.39, -0.9855, 3129.5229),
    @("INDIQUBE", 212.64, -0.7561, 4465.6847),
    @("EBGNG", 346.65, 3.2311, 3952.2092),
    @("LIQGRWBEES", 1014.74, 0.0246, 0),
    @("CHEMBONDCH", 153.35, -1.6987, 412.459),
    @("GROWWNIFTY", 10.29, -0.3872, 0),
    @("ANTHEM", 702.25, -0.1209, 39439.0658),
    @("QUALITY30", 21.05, -0.8945, 0),
    @("SMARTWORKS", 606.65, 2.0867, 6931.2448)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 5).Value = $row[2]
    $ws.Cells.Item($r, 8).Value = $row[3]
}